$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("H111").Value = 774.75
$ws.Range("I111").Value = 849
$ws.Range("K111").Value = 2547
$ws.Range("M111").Value = 520
$ws.Range("H131").Value = 128608.125
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
$ws.Range("H137").Value = 1788.909
$ws.Range("I137").Value = 1521.2307
$ws.Range("J137").Value = 2175.5557
$ws.Range("K137").Value = 4563.6921
$ws.Range("L137").Value = 6526.6671
$ws.Range("M137").Value = -2013.6921
$ws.Range("N137").Value = -11626.6671
$ws.Range("H138").Value = 2243.1875
$ws.Range("I138").Value = 1524.6857
$ws.Range("J138").Value = 3110.3447
$ws.Range("K138").Value = 4574.0571
$ws.Range("L138").Value = 9331.0341
$ws.Range("M138").Value = 565.9429
$ws.Range("N138").Value = -19611.0341

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2603.7468
$ws.Range("I32").Value = 2725.2266
$ws.Range("J32").Value = 326
$ws.Range("K32").Value = 2725.2266
$ws.Range("L32").Value = 326
$ws.Range("M32").Value = -2438.2266
$ws.Range("N32").Value = -900
$ws.Range("H45").Value = 9449.429
$ws.Range("J45").Value = 9981.091
$ws.Range("L45").Value = 9981.091
$ws.Range("N45").Value = -10735.091
$ws.Range("H61").Value = 8193.542
$ws.Range("I61").Value = 6799.8335
$ws.Range("J61").Value = 12374.667
$ws.Range("K61").Value = 6799.8335
$ws.Range("L61").Value = 12374.667
$ws.Range("M61").Value = -6587.8335
$ws.Range("N61").Value = -12798.667
$ws.Range("H63").Value = 1679.5555
$ws.Range("I63").Value = 1679.5555
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 1679.5555
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -993.5554999999999
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 1679.5555
$ws.Range("I66").Value = 1679.5555
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 8397.7775
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -4965.7775
$ws.Range("N66").ClearContents()
$ws.Range("H74").Value = 2740.3225
$ws.Range("I74").Value = 1926.9642
$ws.Range("J74").Value = 10331.667
$ws.Range("K74").Value = 1926.9642
$ws.Range("L74").Value = 10331.667
$ws.Range("M74").Value = -1052.9642
$ws.Range("N74").Value = -12079.667
$ws.Range("H77").Value = 2740.3225
$ws.Range("I77").Value = 1926.9642
$ws.Range("J77").Value = 10331.667
$ws.Range("K77").Value = 9634.821
$ws.Range("L77").Value = 51658.335
$ws.Range("M77").Value = -5266.821
$ws.Range("N77").Value = -60394.335
$ws.Range("H132").Value = 7686.5713
$ws.Range("I132").Value = 6467.6665
$ws.Range("K132").Value = 19402.9995
$ws.Range("M132").Value = -16872.9995
$ws.Range("H136").Value = 8193.542
$ws.Range("I136").Value = 6799.8335
$ws.Range("J136").Value = 12374.667
$ws.Range("K136").Value = 20399.5005
$ws.Range("L136").Value = 37124.001
$ws.Range("M136").Value = -17849.5005
$ws.Range("N136").Value = -42224.001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1746.8
$ws.Range("I105").Value = 1827.8334
$ws.Range("K105").Value = 1827.8334
$ws.Range("M105").Value = -80.83339999999998
$ws.Range("H134").Value = 2085.7322
$ws.Range("I134").Value = 2100.0188
$ws.Range("J134").Value = 1833.3334
$ws.Range("K134").Value = 6300.056399999999
$ws.Range("L134").Value = 5500.0002
$ws.Range("M134").Value = -3765.056399999999
$ws.Range("N134").Value = -10570.0002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 298.30768
$ws.Range("I7").Value = 362.8
$ws.Range("J7").Value = 83.333336
$ws.Range("K7").Value = 362.8
$ws.Range("L7").Value = 83.333336
$ws.Range("M7").Value = -249.8
$ws.Range("N7").Value = -309.333336
$ws.Range("H16").Value = 2924.5881
$ws.Range("I16").Value = 2277.8333
$ws.Range("J16").Value = 4476.8
$ws.Range("K16").Value = 2277.8333
$ws.Range("L16").Value = 4476.8
$ws.Range("M16").Value = -1990.8333
$ws.Range("N16").Value = -5050.8
$ws.Range("H31").Value = 5345.9614
$ws.Range("I31").Value = 3840.875
$ws.Range("K31").Value = 3840.875
$ws.Range("M31").Value = -3545.875
$ws.Range("H34").Value = 5345.9614
$ws.Range("I34").Value = 3840.875
$ws.Range("K34").Value = 3840.875
$ws.Range("M34").Value = -3638.875
$ws.Range("H58").Value = 4179.737
$ws.Range("I58").Value = 1634
$ws.Range("K58").Value = 1634
$ws.Range("M58").Value = -1431
$ws.Range("H94").Value = 452.2
$ws.Range("I94").Value = 287.5
$ws.Range("J94").Value = 1111
$ws.Range("K94").Value = 287.5
$ws.Range("L94").Value = 1111
$ws.Range("M94").Value = 163.5
$ws.Range("N94").Value = -2013
$ws.Range("H95").Value = 14130
$ws.Range("J95").Value = 14130
$ws.Range("L95").Value = 14130
$ws.Range("N95").Value = -19622
$ws.Range("H99").Value = 2813.7144
$ws.Range("I99").Value = 2882.8333
$ws.Range("J99").Value = 2399
$ws.Range("K99").Value = 2882.8333
$ws.Range("L99").Value = 2399
$ws.Range("M99").Value = -1384.8333
$ws.Range("N99").Value = -5395
$ws.Range("H113").Value = 2924.5881
$ws.Range("I113").Value = 2277.8333
$ws.Range("J113").Value = 4476.8
$ws.Range("K113").Value = 2277.8333
$ws.Range("L113").Value = 4476.8
$ws.Range("M113").Value = -107.8332999999998
$ws.Range("N113").Value = -8816.8
$ws.Range("H122").Value = 4702
$ws.Range("I122").Value = 4953.8
$ws.Range("J122").Value = 3862.6667
$ws.Range("K122").Value = 14861.4
$ws.Range("L122").Value = 11588.0001
$ws.Range("M122").Value = -12411.4
$ws.Range("N122").Value = -16488.0001
$ws.Range("H126").Value = 2813.7144
$ws.Range("I126").Value = 2882.8333
$ws.Range("J126").Value = 2399
$ws.Range("K126").Value = 8648.499899999999
$ws.Range("L126").Value = 7197
$ws.Range("M126").Value = -6178.499899999999
$ws.Range("N126").Value = -12137
$ws.Range("H134").Value = 6242.636
$ws.Range("I134").Value = 2958.625
$ws.Range("J134").Value = 15000
$ws.Range("K134").Value = 8875.875
$ws.Range("L134").Value = 45000
$ws.Range("M134").Value = -6340.875
$ws.Range("N134").Value = -50070
$ws.Range("H136").Value = 4179.737
$ws.Range("I136").Value = 1634
$ws.Range("K136").Value = 4902
$ws.Range("M136").Value = -2352

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 997
$ws.Range("I92").Value = 1000
$ws.Range("K92").Value = 3000
$ws.Range("M92").Value = -1752
$ws.Range("H122").Value = 692.8
$ws.Range("I122").Value = 198.5
$ws.Range("J122").Value = 1022.3333
$ws.Range("K122").Value = 1786.5
$ws.Range("L122").Value = 9200.9997
$ws.Range("M122").Value = 663.5
$ws.Range("N122").Value = -14100.9997
$ws.Range("H132").Value = 1629.8182
$ws.Range("I132").Value = 1171.3334
$ws.Range("J132").Value = 2180
$ws.Range("K132").Value = 10542.0006
$ws.Range("L132").Value = 19620
$ws.Range("M132").Value = -8012.000599999999
$ws.Range("N132").Value = -24680

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 6999.8
$ws.Range("I102").Value = 5833
$ws.Range("K102").Value = 5833
$ws.Range("M102").Value = -4211
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5113.2856
$ws.Range("I40").Value = 5113.2856
$ws.Range("K40").Value = 5113.2856
$ws.Range("M40").Value = -4977.2856
$ws.Range("H122").Value = 3234.3572
$ws.Range("J122").Value = 3575.6667
$ws.Range("L122").Value = 10727.0001
$ws.Range("N122").Value = -15627.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 15014.143
$ws.Range("J41").Value = 16619.8
$ws.Range("L41").Value = 16619.8
$ws.Range("N41").Value = -17399.8
$ws.Range("H61").Value = 10000
$ws.Range("J61").Value = 10000
$ws.Range("L61").Value = 10000
$ws.Range("N61").Value = -10584
$ws.Range("H63").Value = 37550.75
$ws.Range("J63").Value = 37550.75
$ws.Range("L63").Value = 37550.75
$ws.Range("N63").Value = -38798.75
$ws.Range("H66").Value = 37550.75
$ws.Range("J66").Value = 37550.75
$ws.Range("L66").Value = 112652.25
$ws.Range("N66").Value = -118892.25
$ws.Range("H82").Value = 49497.5
$ws.Range("J82").Value = 49497.5
$ws.Range("L82").Value = 49497.5
$ws.Range("N82").Value = -50263.5
$ws.Range("H85").Value = 49497.5
$ws.Range("J85").Value = 49497.5
$ws.Range("L85").Value = 49497.5
$ws.Range("N85").Value = -52149.5
$ws.Range("H107").Value = 689.1111
$ws.Range("J107").Value = 760.875
$ws.Range("L107").Value = 2282.625
$ws.Range("N107").Value = -6122.625
$ws.Range("H119").Value = 61747.5
$ws.Range("J119").Value = 61747.5
$ws.Range("L119").Value = 61747.5
$ws.Range("N119").Value = -71423.5
$ws.Range("H136").Value = 4890.8213
$ws.Range("I136").Value = 3234.0625
$ws.Range("J136").Value = 7099.8335
$ws.Range("K136").Value = 9702.1875
$ws.Range("L136").Value = 21299.5005
$ws.Range("M136").Value = -7152.1875
$ws.Range("N136").Value = -26399.5005
